# Update the "Apps" column (Z) on the "devices" sheet:
#  - Merge "Podcast" and "Musics" into a single "Audio" token (keeping the
#    position of the first occurrence, dropping any duplicate).
#  - Replace "AllTrails" with "Maps" unless "Maps" is already present in the
#    list, in which case "AllTrails" is simply dropped.
#  - Normalize the separator from "," to ", " (comma + space).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("devices")

$lastRow = 384
$col = 26  # column Z

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $raw = $cell.Value2

    if ($null -eq $raw) { continue }
    $original = [string]$raw
    if ($original -eq "") { continue }

    $tokens = $original.Split(",")

    $mapsPresent = $false
    foreach ($tok in $tokens) {
        if ($tok -eq "Maps") { $mapsPresent = $true }
    }

    $audioAdded = $false
    $mapsAdded = $false
    $newValue = ""
    $first = $true

    foreach ($tok in $tokens) {
        $piece = $null
        if (($tok -eq "Podcast") -or ($tok -eq "Musics")) {
            if (-not $audioAdded) {
                $piece = "Audio"
                $audioAdded = $true
            }
        }
        elseif ($tok -eq "AllTrails") {
            if (-not $mapsPresent) {
                if (-not $mapsAdded) {
                    $piece = "Maps"
                    $mapsAdded = $true
                }
            }
        }
        else {
            $piece = $tok
        }

        if ($null -ne $piece) {
            if ($first) {
                $newValue = $piece
                $first = $false
            }
            else {
                $newValue = $newValue + ", " + $piece
            }
        }
    }

    if ($newValue -ne $original) {
        $cell.Value = $newValue
    }
}
